$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.345.62"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.274.12"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "308.85"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "97.51"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("D7").Value = "0.529"
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").Value = "35.24"
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("D11").Value = "0.0807"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "6.84"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "2.625.12"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "14.64"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "2.263.01"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "0.792"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "42.201.13"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "12.27"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").Value = "5.99"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "67.67"
$ws.Range("D23").Value = "236.86"
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("D24").Value = "2.59"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "23.60"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "37.38"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("D29").Value = "9.58"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "163.82"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").Value = "5.26"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("D35").Value = "17.71"
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("D36").Value = "0.0736"
$ws.Range("E36").Value = "  -2.47%  "
$ws.Range("D37").Value = "2.38"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").Value = "  -3.91%  "
$ws.Range("D43").Value = "1.948.01"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").Value = "18.79"
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").Value = "9.80"
$ws.Range("E47").Value = "  -4.05%  "
$ws.Range("D48").Value = "53.94"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "2.496.91"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "92.34"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "71.58"
$ws.Range("E51").Value = "  -2.02%  "

# Row 39/40 swap
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "1.82"
$ws.Range("E39").Value = "  -3.39%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.115"
$ws.Range("E40").Value = "  -0.95%  "
